$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.37000000000053
$ws.Range("H2").Value = 0.0000000001378129832474428
$ws.Range("I2").Value = 0.0000000001378129832474428
$ws.Range("L2").Value = 57.73507357332018
$ws.Range("M2").Value = "[43.597460256418174, 71.87268689022218]"
$ws.Range("N2").Value = 0.0000000001624487211415726
$ws.Range("O2").Value = 0.0000000001624487211415726
$ws.Range("P2").Value = 1.578658170272348
$ws.Range("Q2").Value = "[1.2893423303021159, 1.8679740102425804]"
$ws.Range("R2").Value = 0.00000000000002486899575160351
$ws.Range("S2").Value = 0.00000000000002486899575160351
$ws.Range("T2").Value = 54.28129269769806
$ws.Range("U2").Value = "[44.746816287618216, 63.81576910777791]"
$ws.Range("V2").Value = 0.000000000000005995204332975845
$ws.Range("W2").Value = 0.000000000000005995204332975845
$ws.Range("X2").Value = 18.99575575575615
$ws.Range("Y2").Value = 17.82756756756793
$ws.Range("Z2").Value = 20.16394394394437
$ws.Range("F3").Value = 25.37000000000053
$ws.Range("H3").Value = 0.00000000000005062616992290714
$ws.Range("I3").Value = 0.00000000000005062616992290714
$ws.Range("L3").Value = 59.32042885740093
$ws.Range("M3").Value = "[46.28726237136905, 72.35359534343281]"
$ws.Range("N3").Value = 0.000000000007368994303647014
$ws.Range("O3").Value = 0.000000000007368994303647014
$ws.Range("P3").Value = 1.239026532046425
$ws.Range("Q3").Value = "[1.0000264903318863, 1.478026573760963]"
$ws.Range("R3").Value = 0.0000000000001318944953254686
$ws.Range("S3").Value = 0.0000000000001318944953254686
$ws.Range("T3").Value = 50.68093369882483
$ws.Range("U3").Value = "[43.17999582223792, 58.18187157541174]"
$ws.Range("X3").Value = 20.36710710710753
$ws.Range("Y3").Value = 19.40208208208249
$ws.Range("Z3").Value = 21.33213213213257
$ws.Range("F4").Value = 25.37000000000053
$ws.Range("H4").Value = 0.00000000000002964295475749168
$ws.Range("I4").Value = 0.00000000000002964295475749168
$ws.Range("L4").Value = 59.56693893621589
$ws.Range("M4").Value = "[45.44691295923824, 73.68696491319355]"
$ws.Range("N4").Value = 0.00000000006598077639807798
$ws.Range("O4").Value = 0.00000000006598077639807798
$ws.Range("P4").Value = 0.9371317425122703
$ws.Range("Q4").Value = "[0.7107106503616549, 1.1635528346628856]"
$ws.Range("R4").Value = 0.0000000001123143800185744
$ws.Range("S4").Value = 0.0000000001123143800185744
$ws.Range("T4").Value = 55.69367738826798
$ws.Range("U4").Value = "[48.406208580995035, 62.98114619554092]"
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 21.58608608608654
$ws.Range("Y4").Value = 20.67185185185228
$ws.Range("Z4").Value = 22.50032032032079
$ws.Range("F5").Value = 25.37000000000053
$ws.Range("H5").Value = 0.000000001073317767996684
$ws.Range("I5").Value = 0.000000001073317767996684
$ws.Range("L5").Value = 54.21100523322923
$ws.Range("M5").Value = "[36.90111363539084, 71.52089683106762]"
$ws.Range("N5").Value = 0.0000001089146637323068
$ws.Range("O5").Value = 0.0000001089146637323068
$ws.Range("P5").Value = 0.5094474573388856
$ws.Range("Q5").Value = "[0.1823947686768852, 0.8365001460008861]"
$ws.Range("R5").Value = 0.00300437383607588
$ws.Range("S5").Value = 0.00300437383607588
$ws.Range("T5").Value = 58.52152969844411
$ws.Range("U5").Value = "[49.36974718924026, 67.67331220764795]"
$ws.Range("V5").Value = 0.0000000000000002220446049250313
$ws.Range("W5").Value = 0.0000000000000002220446049250313
$ws.Range("X5").Value = 23.31297297297345
$ws.Range("Y5").Value = 21.99241241241287
$ws.Range("Z5").Value = 24.63353353353404
$ws.Range("F6").Value = 22.94000000000015
$ws.Range("H6").Value = 0.000000000003593569886106707
$ws.Range("I6").Value = 0.000000000003593569886106707
$ws.Range("L6").Value = 63.67814158511832
$ws.Range("M6").Value = "[49.63807087110929, 77.71821229912734]"
$ws.Range("N6").Value = 0.000000000008180789379252928
$ws.Range("O6").Value = 0.000000000008180789379252928
$ws.Range("P6").Value = -0.03773684869176908
$ws.Range("Q6").Value = "[-0.28931583997023136, 0.2138421425866932]"
$ws.Range("R6").Value = 0.7639559724542799
$ws.Range("S6").Value = 0.7639559724542799
$ws.Range("T6").Value = 55.78020463825701
$ws.Range("U6").Value = "[46.52591975987619, 65.03448951663783]"
$ws.Range("V6").Value = 0.0000000000000008881784197001252
$ws.Range("W6").Value = 0.0000000000000008881784197001252
$ws.Range("X6").Value = 0.1377777777777816
$ws.Range("Y6").Value = -0.7807407407407443
$ws.Range("Z6").Value = 1.056296296296307
$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 22.94000000000015
$ws.Range("H7").Value = 0.0000000001156768014709542
$ws.Range("I7").Value = 0.0000000001156768014709542
$ws.Range("L7").Value = 53.58559354557613
$ws.Range("M7").Value = "[38.62331752062666, 68.5478695705256]"
$ws.Range("N7").Value = 0.000000004924323393495911
$ws.Range("O7").Value = 0.000000004924323393495911
$ws.Range("P7").Value = 0.2327105669325773
$ws.Range("Q7").Value = "[-0.08176317216550011, 0.5471843060306547]"
$ws.Range("R7").Value = 0.1430861536172205
$ws.Range("S7").Value = 0.1430861536172205
$ws.Range("T7").Value = 50.60863025496415
$ws.Range("U7").Value = "[41.89290110155701, 59.3243594083713]"
$ws.Range("V7").Value = 0.000000000000003108624468950438
$ws.Range("W7").Value = 0.000000000000003108624468950438
$ws.Range("X7").Value = 22.09037037037051
$ws.Range("Y7").Value = 20.94222222222236
$ws.Range("Z7").Value = 23.23851851851867
$ws.Range("F8").Value = 22.94000000000015
$ws.Range("H8").Value = 0.0000000000006020739462542224
$ws.Range("I8").Value = 0.0000000000006020739462542224
$ws.Range("L8").Value = 59.67992353629768
$ws.Range("M8").Value = "[44.6589920768879, 74.70085499570746]"
$ws.Range("N8").Value = 0.0000000003420956851130086
$ws.Range("O8").Value = 0.0000000003420956851130086
$ws.Range("P8").Value = 0.4842895582110396
$ws.Range("Q8").Value = "[0.22013161736865516, 0.748447499053424]"
$ws.Range("R8").Value = 0.0005983081947840141
$ws.Range("S8").Value = 0.0005983081947840141
$ws.Range("T8").Value = 55.87995031025161
$ws.Range("U8").Value = "[47.77696305585222, 63.982937564650996]"
$ws.Range("X8").Value = 21.17185185185199
$ws.Range("Y8").Value = 20.20740740740754
$ws.Range("Z8").Value = 22.13629629629644
